# Map119.xlsx "huge v2 update" translation patch
# - Inserts a new shared string "ーーーーー基本変更点ーーーーー" as a new row 9
#   (pushing the former rows 9..98 down to 10..99).
# - Collapses the old 4-column layout (A = source text, C/D = translated
#   text, duplicated) into a 2-column layout (A = source text, B =
#   translated text). Where no translation existed (no D value), B mirrors
#   A (self-referential placeholder for "not yet translated").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new row for the section marker, shifting everything at/after
#    row 9 down by one. This reproduces the net effect of the new shared
#    string being spliced into the middle of the table.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value2 = "ーーーーー基本変更点ーーーーー"

# 2) Rebuild column B from column D (or mirror column A when there was no
#    D translation), then drop the now-unused C/D columns.
for ($r = 1; $r -le 99; $r++) {
    if ($r -ne 9) {
        $aVal = $ws.Cells.Item($r, 1).Value2
        $dVal = $ws.Cells.Item($r, 4).Value2
        if ($dVal -ne $null -and $dVal -ne "") {
            $ws.Cells.Item($r, 2).Value2 = $dVal
        } else {
            $ws.Cells.Item($r, 2).Value2 = $aVal
        }
    }
    $ws.Cells.Item($r, 3).ClearContents()
    $ws.Cells.Item($r, 4).ClearContents()
}
